{"js": "// Apply the dated worksheet update: refresh the header date and all\n// two-digit multiplication problems/answers in the table.\nconst replacements = [\n  [\"2025-09-11 Thursday\", \"2025-09-12 Friday\"],\n  [\"22\u00d761=1342\", \"11\u00d763=693\"],\n  [\"71\u00d775=5325\", \"84\u00d749=4116\"],\n  [\"76\u00d799=7524\", \"11\u00d711=121\"],\n  [\"74\u00d774=5476\", \"32\u00d799=3168\"],\n  [\"68\u00d777=5236\", \"21\u00d742=882\"],\n  [\"56\u00d718=1008\", \"82\u00d738=3116\"],\n  [\"84\u00d726=2184\", \"52\u00d721=1092\"],\n  [\"65\u00d745=2925\", \"44\u00d780=3520\"],\n  [\"78\u00d717=1326\", \"51\u00d764=3264\"],\n  [\"98\u00d781=7938\", \"70\u00d735=2450\"],\n  [\"77\u00d792=7084\", \"21\u00d794=1974\"],\n  [\"14\u00d721=294\", \"91\u00d714=1274\"],\n  [\"82\u00d757=4674\", \"58\u00d797=5626\"],\n  [\"82\u00d789=7298\", \"78\u00d759=4602\"],\n  [\"32\u00d741=1312\", \"25\u00d732=800\"],\n  [\"96\u00d713=1248\", \"18\u00d726=468\"],\n  [\"61\u00d759=3599\", \"63\u00d741=2583\"],\n  [\"83\u00d781=6723\", \"39\u00d761=2379\"],\n  [\"26\u00d767=1742\", \"59\u00d773=4307\"],\n  [\"87\u00d726=2262\", \"85\u00d746=3910\"],\n  [\"85\u00d726=2210\", \"36\u00d764=2304\"],\n  [\"31\u00d751=1581\", \"63\u00d713=819\"],\n  [\"36\u00d791=3276\", \"99\u00d712=1188\"],\n  [\"81\u00d770=5670\", \"65\u00d751=3315\"],\n  [\"45\u00d771=3195\", \"21\u00d740=840\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n}\n\nawait context.sync();\n\nreturn `replaced ${totalReplaced} of ${replacements.length}`;", "ps1": "# Apply the dated worksheet update: refresh the header date and all\n# two-digit multiplication problems/answers in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-11 Thursday\", \"2025-09-12 Friday\"),\n    @(\"22\u00d761=1342\", \"11\u00d763=693\"),\n    @(\"71\u00d775=5325\", \"84\u00d749=4116\"),\n    @(\"76\u00d799=7524\", \"11\u00d711=121\"),\n    @(\"74\u00d774=5476\", \"32\u00d799=3168\"),\n    @(\"68\u00d777=5236\", \"21\u00d742=882\"),\n    @(\"56\u00d718=1008\", \"82\u00d738=3116\"),\n    @(\"84\u00d726=2184\", \"52\u00d721=1092\"),\n    @(\"65\u00d745=2925\", \"44\u00d780=3520\"),\n    @(\"78\u00d717=1326\", \"51\u00d764=3264\"),\n    @(\"98\u00d781=7938\", \"70\u00d735=2450\"),\n    @(\"77\u00d792=7084\", \"21\u00d794=1974\"),\n    @(\"14\u00d721=294\", \"91\u00d714=1274\"),\n    @(\"82\u00d757=4674\", \"58\u00d797=5626\"),\n    @(\"82\u00d789=7298\", \"78\u00d759=4602\"),\n    @(\"32\u00d741=1312\", \"25\u00d732=800\"),\n    @(\"96\u00d713=1248\", \"18\u00d726=468\"),\n    @(\"61\u00d759=3599\", \"63\u00d741=2583\"),\n    @(\"83\u00d781=6723\", \"39\u00d761=2379\"),\n    @(\"26\u00d767=1742\", \"59\u00d773=4307\"),\n    @(\"87\u00d726=2262\", \"85\u00d746=3910\"),\n    @(\"85\u00d726=2210\", \"36\u00d764=2304\"),\n    @(\"31\u00d751=1581\", \"63\u00d713=819\"),\n    @(\"36\u00d791=3276\", \"99\u00d712=1188\"),\n    @(\"81\u00d770=5670\", \"65\u00d751=3315\"),\n    @(\"45\u00d771=3195\", \"21\u00d740=840\"),\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\nWrite-Output \"done\""}
